# Update churn-rate percentages across the three sheets to reflect
# the refreshed source data (values recomputed upstream).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Без группировок")
$ws2 = $wb.Worksheets.Item("По странам")
$ws3 = $wb.Worksheets.Item("По типу подписки")

# Без группировок
$ws1.Range("B2").Value = 0.953968253968254
$ws1.Range("C2").Value = 0.9408730158730159
$ws1.Range("D2").Value = 0.9293650793650794
$ws1.Range("E2").Value = 0.8214285714285714
$ws1.Range("F2").Value = 0.8067460317460318
$ws1.Range("G2").Value = 0.7936507936507936
$ws1.Range("H2").Value = 0.7813492063492063
$ws1.Range("I2").Value = 0.7753968253968254
$ws1.Range("J2").Value = 0.7686507936507937
$ws1.Range("K2").Value = 0.7623015873015873
$ws1.Range("L2").Value = 0.7607142857142857
$ws1.Range("M2").Value = 0.7634920634920634
$ws1.Range("N2").Value = 0.7523809523809524
$ws1.Range("O2").Value = 0.7154761904761905
$ws1.Range("P2").Value = 0.6869047619047619
$ws1.Range("Q2").Value = 0.5595238095238095

# По странам
$ws2.Range("D2").Value = 0.1444444444444444
$ws2.Range("E2").Value = 0.1894736842105263
$ws2.Range("F2").Value = 0.1891891891891892
$ws2.Range("G2").Value = 0.009523809523809525
$ws2.Range("H2").Value = 0.02608695652173913
$ws2.Range("I2").Value = 0.04225352112676056
$ws2.Range("J2").Value = 0.0131578947368421
$ws2.Range("K2").Value = 0.04320987654320987
$ws2.Range("L2").Value = 0.05142857142857143
$ws2.Range("M2").Value = 0.1183431952662722
$ws2.Range("N2").Value = 0.06918238993710692
$ws2.Range("O2").Value = 0.07428571428571429
$ws2.Range("P2").Value = 0.05699481865284974
$ws2.Range("Q2").Value = 0.08465608465608465
$ws2.Range("R2").Value = 0.09852216748768473
$ws2.Range("S2").Value = 0.08900523560209424
$ws2.Range("D3").Value = 0.1428571428571428
$ws2.Range("G3").Value = 0
$ws2.Range("H3").Value = 0
$ws2.Range("L3").Value = 0.09090909090909091
$ws2.Range("M3").Value = 0.02439024390243903
$ws2.Range("N3").Value = 0
$ws2.Range("P3").Value = 0.1132075471698113
$ws2.Range("Q3").Value = 0.07843137254901961
$ws2.Range("R3").Value = 0.1818181818181818
$ws2.Range("D4").Value = 0
$ws2.Range("E4").Value = 0.01219512195121951
$ws2.Range("F4").Value = 0.03448275862068965
$ws2.Range("G4").Value = 0.03523035230352303
$ws2.Range("H4").Value = 0.03341902313624678
$ws2.Range("I4").Value = 0.05527638190954774
$ws2.Range("J4").Value = 0.06483790523690773
$ws2.Range("K4").Value = 0.03508771929824561
$ws2.Range("L4").Value = 0.06205250596658711
$ws2.Range("M4").Value = 0.05263157894736842
$ws2.Range("N4").Value = 0.08558558558558559
$ws2.Range("O4").Value = 0.06487695749440715
$ws2.Range("P4").Value = 0.07400379506641366
$ws2.Range("Q4").Value = 0.08667736757624397
$ws2.Range("R4").Value = 0.1073619631901841
$ws2.Range("S4").Value = 0.07918781725888324

# По типу подписки
$ws3.Range("D2").Value = 0.1428571428571428
$ws3.Range("E2").Value = 0.15
$ws3.Range("F2").Value = 0.08644859813084112
$ws3.Range("G2").Value = 0.03167420814479638
$ws3.Range("H2").Value = 0.03354297693920336
$ws3.Range("I2").Value = 0.05384615384615385
$ws3.Range("J2").Value = 0.05400372439478585
$ws3.Range("K2").Value = 0.04638218923933209
$ws3.Range("L2").Value = 0.06177606177606178
$ws3.Range("M2").Value = 0.0653061224489796
$ws3.Range("N2").Value = 0.0670995670995671
$ws3.Range("O2").Value = 0.05707762557077625
$ws3.Range("P2").Value = 0.06038647342995169
$ws3.Range("Q2").Value = 0.0586734693877551
$ws3.Range("S2").Value = 0.03846153846153846
$ws3.Range("J3").Value = 0
$ws3.Range("K3").Value = 0.01428571428571429
$ws3.Range("L3").Value = 0.05833333333333333
$ws3.Range("M3").Value = 0.07643312101910828
$ws3.Range("N3").Value = 0.09836065573770492
$ws3.Range("O3").Value = 0.08260869565217391
$ws3.Range("P3").Value = 0.08635097493036212
$ws3.Range("Q3").Value = 0.1082802547770701
$ws3.Range("R3").Value = 0.1189931350114417
$ws3.Range("S3").Value = 0.09354120267260579

